$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("research")

$NL = [char]10

# --- Insert a new row at position 3, shifting the old row 3 ("10fold_1M_a") down to row 4 ---
$ws.Rows.Item(3).Insert()

# --- Give the new row 3 the same thin-border style used by the rest of the table ---
$ws.Range("A3:E3").Borders.LineStyle = 1

# --- Populate the DESCRIPTION cells first (this order reproduces the original shared-string layout) ---
$ws.Range("E3").Value = "10 folds, each having 10000 playlists. " + $NL + "Includes categories 1, 2, 3."
$ws.Range("E2").Value = "10 folds, each having 10000 playlists. " + $NL + "Includes only category 1."
$ws.Range("E4").Value = "10 folds, each having 100000 playlists. " + $NL + "Includes only category 1."

# --- Populate the rest of the new row 3 ("10fold_100K_b") ---
$ws.Range("A3").Value = "10fold_100K_b"
$ws.Range("B3").Value = "100K"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 3

# --- Wrap the text of the DESCRIPTION column for the data rows and size the rows to fit two lines ---
$ws.Range("E2").WrapText = $true
$ws.Range("E3").WrapText = $true
$ws.Range("E4").WrapText = $true

$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30

# --- Adjust column widths to fit the new, longer content ---
$ws.Columns.Item(1).ColumnWidth = 19.14
$ws.Columns.Item(2).ColumnWidth = 14.14
$ws.Columns.Item(5).ColumnWidth = 70.57

# --- Put the selection on the new last row, as in the saved file ---
$ws.Range("A4").Select()
